$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole table one column to the right (A->B, B->C, C->D, D->E),
# bringing along values and formatting (xlShiftToRight = -4161).
$ws.Range("A1:A15").Insert(-4161)

# The new column A (rows 2-15) should carry the same header-row style (s="1")
# as the rest of the table, but stay empty. Copy formatting only from B1.
$ws.Range("B1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)

# Tweak label text: add spaces around "=" in the n= and P= annotations.
$ws.Range("D1").Value = "Control at T1 (n = 745)"
$ws.Range("E1").Value = "Control at T2 (n = 745)"
$ws.Range("B3").Value = "Gender (P = 1.000)"
$ws.Range("B10").Value = "Interested in News (P = 1.000)"
